# Updated cryptos list — apply Price (D) and Volume(1h) (E) changes per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.832.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.379.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.376.19"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.79%  "

$ws.Range("E9").Value = "  -1.14%  "

$ws.Range("E10").Value = "  +2.36%  "

$ws.Range("E11").Value = "  -2.32%  "

$ws.Range("E12").Value = "  -1.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.959.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.66%  "

$ws.Range("E14").Value = "  +0.38%  "

$ws.Range("E15").Value = "  +1.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.390.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.039.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.85%  "

$ws.Range("E20").Value = "  -1.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "374.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.61%  "

$ws.Range("E23").Value = "  -2.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.527.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("E26").Value = "  -2.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.55%  "

$ws.Range("E28").Value = "  +11.32%  "

$ws.Range("E29").Value = "  +8.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.990"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.02%  "

$ws.Range("E33").Value = "  -1.38%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.21%  "

$ws.Range("E37").Value = "  -1.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.87%  "

$ws.Range("E40").Value = "  -3.51%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.777"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("E46").Value = "  -3.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.452.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.33%  "

$ws.Range("E50").Value = "  -2.46%  "

$ws.Range("E51").Value = "  +4.62%  "
